$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.637.67"
$ws.Range("E2").Value = "  -1.34%  "
$ws.Range("D3").Value = "2.522.99"
$ws.Range("E3").Value = "  +3.02%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.81%  "
$ws.Range("D5").Value = "309.46"
$ws.Range("E5").Value = "  -2.02%  "
$ws.Range("D6").Value = "100.50"
$ws.Range("E6").Value = "  +2.33%  "
$ws.Range("D7").Value = "0.569"
$ws.Range("E7").Value = "  -1.70%  "
$ws.Range("E8").Value = "  +1.16%  "
$ws.Range("D9").Value = "0.525"
$ws.Range("E9").Value = "  -2.91%  "
$ws.Range("D10").Value = "35.85"
$ws.Range("E10").Value = "  -0.48%  "
$ws.Range("D11").Value = "0.0802"
$ws.Range("E11").Value = "  -1.63%  "
$ws.Range("D12").Value = "7.27"
$ws.Range("E12").Value = "  -3.76%  "
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("D14").Value = "2.911.42"
$ws.Range("E14").Value = "  +2.27%  "
$ws.Range("D15").Value = "15.62"
$ws.Range("E15").Value = "  +2.19%  "
$ws.Range("D16").Value = "2.545.74"
$ws.Range("E16").Value = "  -1.94%  "
$ws.Range("D17").Value = "0.804"
$ws.Range("E17").Value = "  -5.53%  "
$ws.Range("D18").Value = "42.613.01"
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("D19").Value = "6.71"
$ws.Range("D20").Value = "0.0₃0948"
$ws.Range("D21").Value = "12.13"
$ws.Range("E21").Value = "  -5.49%  "
$ws.Range("D22").Value = "69.28"
$ws.Range("E22").Value = "  -0.99%  "
$ws.Range("D23").Value = "243.88"
$ws.Range("E23").Value = "  -4.41%  "
$ws.Range("D24").Value = "2.88"
$ws.Range("E24").Value = "  -3.05%  "
$ws.Range("E25").Value = "  -3.17%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "26.18"
$ws.Range("E27").Value = "  -4.27%  "
$ws.Range("E28").Value = "  -3.87%  "
$ws.Range("D29").Value = "38.98"
$ws.Range("E29").Value = "  -5.36%  "
$ws.Range("D30").Value = "10.13"
$ws.Range("E30").Value = "  -2.22%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "155.82"
$ws.Range("E31").Value = "  -0.52%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "5.76"
$ws.Range("E32").Value = "  -2.30%  "
$ws.Range("E33").Value = "  +10.24%  "
$ws.Range("D34").Value = "0.0785"
$ws.Range("E34").Value = "  -3.52%  "
$ws.Range("E35").Value = "  -2.87%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "2.03"
$ws.Range("E36").Value = "  -6.39%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "3.18"
$ws.Range("E37").Value = "  -8.07%  "
$ws.Range("D38").Value = "18.14"
$ws.Range("E38").Value = "  -4.08%  "
$ws.Range("E39").Value = "  -1.54%  "
$ws.Range("E40").Value = "  -0.49%  "
$ws.Range("D41").Value = "4.22"
$ws.Range("E41").Value = "  +4.99%  "
$ws.Range("D42").Value = "22.04"
$ws.Range("E42").Value = "  -5.00%  "
$ws.Range("E43").Value = "  +0.38%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "0.0299"
$ws.Range("E44").Value = "  -2.16%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "3.27"
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("D46").Value = "1.988.03"
$ws.Range("E46").Value = "  -1.27%  "
$ws.Range("D47").Value = "8.81"
$ws.Range("E47").Value = "  -2.19%  "
$ws.Range("D48").Value = "2.766.63"
$ws.Range("E48").Value = "  -1.50%  "
$ws.Range("D49").Value = "80.05"
$ws.Range("E49").Value = "  -4.10%  "
$ws.Range("E50").Value = "  -3.53%  "
$ws.Range("D51").Value = "72.24"
$ws.Range("E51").Value = "  -3.45%  "
